# Updates the cryptos list (Price / Volume(1h) columns) to refreshed values.
# Note: some new Price values look like plain numbers (e.g. "0.998"); those
# are written with a leading apostrophe so Excel stores them as text (as the
# source data already is) instead of auto-converting them to numeric cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '76.341.79'
$ws.Range("E2").Value = '  +0.32%  '
$ws.Range("D3").Value = '3.041.88'
$ws.Range("E3").Value = '  +4.05%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = "'198.57"
$ws.Range("E5").Value = '  -0.52%  '
$ws.Range("D6").Value = "'617.01"
$ws.Range("E6").Value = '  +3.47%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("E8").Value = '  -0.50%  '
$ws.Range("D9").Value = "'0.205"
$ws.Range("E9").Value = '  +4.40%  '
$ws.Range("D10").Value = '3.042.60'
$ws.Range("E10").Value = '  +3.63%  '
$ws.Range("E11").Value = '  -1.68%  '
$ws.Range("E12").Value = '  -0.43%  '
$ws.Range("E13").Value = '  +6.13%  '
$ws.Range("D14").Value = '3.601.16'
$ws.Range("E14").Value = '  +3.60%  '
$ws.Range("D15").Value = "'28.81"
$ws.Range("E15").Value = '  +3.18%  '
$ws.Range("D16").Value = '76.264.63'
$ws.Range("E16").Value = '  +0.47%  '
$ws.Range("E17").Value = '  +2.03%  '
$ws.Range("D18").Value = '3.047.62'
$ws.Range("E18").Value = '  +4.57%  '
$ws.Range("D19").Value = "'13.52"
$ws.Range("E19").Value = '  +1.74%  '
$ws.Range("D20").Value = "'8.96"
$ws.Range("E20").Value = '  +2.52%  '
$ws.Range("D21").Value = "'381.36"
$ws.Range("E21").Value = '  +2.57%  '
$ws.Range("D22").Value = "'2.38"
$ws.Range("E22").Value = '  +3.58%  '
$ws.Range("E23").Value = '  +1.20%  '
$ws.Range("D24").Value = '3.199.69'
$ws.Range("E24").Value = '  +4.64%  '
$ws.Range("D25").Value = "'72.47"
$ws.Range("E25").Value = '  -0.02%  '
$ws.Range("E26").Value = '  +0.18%  '
$ws.Range("E27").Value = '  +1.87%  '
$ws.Range("D28").Value = "'9.75"
$ws.Range("E28").Value = '  +1.18%  '
$ws.Range("E29").Value = '  +0.31%  '
$ws.Range("D30").Value = "'0.998"
$ws.Range("E30").Value = '  -0.46%  '
$ws.Range("E31").Value = '  +5.01%  '
$ws.Range("D32").Value = "'1.39"
$ws.Range("E32").Value = '  +1.16%  '
$ws.Range("D33").Value = "'493.12"
$ws.Range("E33").Value = '  -0.79%  '
$ws.Range("D34").Value = "'1.92"
$ws.Range("E34").Value = '  +4.62%  '
$ws.Range("E35").Value = '  +0.05%  '
$ws.Range("D36").Value = "'20.54"
$ws.Range("E36").Value = '  +2.21%  '
$ws.Range("D37").Value = "'163.12"
$ws.Range("E37").Value = '  -0.81%  '
$ws.Range("E38").Value = '  +6.29%  '
$ws.Range("E39").Value = '  +1.90%  '
$ws.Range("D40").Value = "'191.51"
$ws.Range("E40").Value = '  +7.40%  '
$ws.Range("E41").Value = '  -2.53%  '
$ws.Range("E42").Value = '  -4.99%  '
$ws.Range("E43").Value = '  +0.02%  '
$ws.Range("D44").Value = "'0.798"
$ws.Range("E44").Value = '  +21.48%  '
$ws.Range("D45").Value = "'5.10"
$ws.Range("E45").Value = '  +3.64%  '
$ws.Range("D46").Value = "'41.93"
$ws.Range("E46").Value = '  +4.53%  '
$ws.Range("E47").Value = '  +4.44%  '
$ws.Range("E48").Value = '  -0.37%  '
$ws.Range("E49").Value = '  +5.39%  '
$ws.Range("D50").Value = "'0.599"
$ws.Range("E50").Value = '  +2.94%  '
$ws.Range("D51").Value = "'3.87"
$ws.Range("E51").Value = '  -0.01%  '
